$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add two new fuzzy rules for TONE_VOICE after the existing rows
$ws.Range("A15").Value = "IF (ROLE_SUBORDINATE IS low) THEN (TONE_VOICE IS mid_tone)"
$ws.Range("B15").Value = "chat_qualifier"

$ws.Range("A16").Value = "IF (ROLE_SUBORDINATE IS high) THEN (TONE_VOICE IS low_tone)"
$ws.Range("B16").Value = "chat_qualifier"

# Fix the PERSONAL_VOCABULARY rule text (vhpc -> vvhpc)
$ws.Range("A7").Value = "IF (PERSONAL_VOCABULARY IS high) THEN (PERSONAL IS vvhpc)"

# Update the selected cell to match the saved view state
$ws.Range("A8").Select()
